$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ivab"
$ws.Range("A3").Value = "nikod"
$ws.Range("A4").Value = "loris"

$ws.Range("C2").Value = "antons"
$ws.Range("C3").Value = "anilf"
$ws.Range("C4").Value = "hoyw"

$ws.Range("E2").Value = "ivab123"
$ws.Range("E3").Value = "nikod321"
$ws.Range("E4").Value = "lorisl321"

$wb.Save()
